$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.959.24'
$ws.Range("E2").Value = '  -2.65%  '
$ws.Range("D3").Value = '3.356.04'
$ws.Range("E3").Value = '  -2.38%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = "'566.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.00%  '
$ws.Range("D6").Value = "'147.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.09%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  +0.35%  '
$ws.Range("E9").Value = '  -0.47%  '
$ws.Range("E10").Value = '  -1.04%  '
$ws.Range("D11").Value = "'0.416"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.19%  '
$ws.Range("D12").Value = '3.932.90'
$ws.Range("E12").Value = '  -2.30%  '
$ws.Range("E13").Value = '  +0.48%  '
$ws.Range("E14").Value = '  -0.59%  '
$ws.Range("D15").Value = '3.359.34'
$ws.Range("E15").Value = '  -2.16%  '
$ws.Range("E16").Value = '  -1.19%  '
$ws.Range("D17").Value = '60.995.39'
$ws.Range("E17").Value = '  -2.66%  '
$ws.Range("E18").Value = '  -0.99%  '
$ws.Range("E19").Value = '  -0.87%  '
$ws.Range("E20").Value = '  -1.43%  '
$ws.Range("D21").Value = "'376.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.83%  '
$ws.Range("D22").Value = "'0.560"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.18%  '
$ws.Range("D23").Value = "'74.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.26%  '
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("E25").Value = '  -2.26%  '
$ws.Range("E26").Value = '  -6.04%  '
$ws.Range("E27").Value = '  -3.68%  '
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.41%  '
$ws.Range("E29").Value = '  -2.40%  '
$ws.Range("E30").Value = '  -0.01%  '
$ws.Range("E31").Value = '  -0.94%  '
$ws.Range("E32").Value = '  -3.58%  '
$ws.Range("D33").Value = "'22.87"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.38%  '
$ws.Range("E34").Value = '  -2.59%  '
$ws.Range("E35").Value = '  +0.72%  '
$ws.Range("D36").Value = "'169.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.44%  '
$ws.Range("E37").Value = '  -3.64%  '
$ws.Range("D38").Value = "'6.81"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.22%  '
$ws.Range("E39").Value = '  -9.82%  '
$ws.Range("D40").Value = '3.392.06'
$ws.Range("E41").Value = '  -2.90%  '
$ws.Range("E42").Value = '  -3.41%  '
$ws.Range("E43").Value = '  -0.82%  '
$ws.Range("E44").Value = '  -5.06%  '
$ws.Range("E45").Value = '  -3.33%  '
$ws.Range("D46").Value = '2.493.66'
$ws.Range("E46").Value = '  -2.41%  '
$ws.Range("D47").Value = "'22.66"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.20%  '
$ws.Range("D48").Value = "'6.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.98%  '
$ws.Range("E49").Value = '  -0.04%  '
$ws.Range("E50").Value = '  -1.85%  '
$ws.Range("D51").Value = "'0.814"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.35%  '
